$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17 (pushes existing rows 17-27 down to 18-28)
$ws.Rows.Item(17).EntireRow.Insert()

# Populate the new row 17 with the latest weekly price observation,
# matching the static columns used throughout this sheet for this series.
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44966
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 100112017
$ws.Cells.Item(17, 7).Value = "Ramas de apio"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 40
$ws.Cells.Item(17, 11).Value = 5000
$ws.Cells.Item(17, 12).Value = 5000
$ws.Cells.Item(17, 13).Value = 5000
$ws.Cells.Item(17, 14).Value = "$/paquete"
$ws.Cells.Item(17, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(17, 16).Value = 5000
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = "Hortaliza"
